# This deck embeds a PowerPoll Office Add-in (a "we:webextension" task-pane
# add-in) on slide 1, rendered via an <mc:AlternateContent> block: a live
# <p:graphicFrame>/<we:webextensionref> for PowerPoint, with a static
# <p:pic> fallback snapshot for hosts that can't render add-ins.
#
# The corresponding commit ("added notifications for incorrect input,
# changed the refresh rate of the graphs") edits the add-in's own
# JavaScript/HTML (outside this .pptx). The only change inside the .pptx
# itself is PowerPoint re-touching the add-in instance when the deck was
# re-saved: every relationship id in ppt/presentation.xml was
# re-minted (sldMasterId/sldId/sldLayoutId/webextensionref/blip r:ids) and
# ppt/slides/udata/data.xml's <we:webextension id="..."/> GUID was
# reassigned. None of the slide's visible content, shapes, geometry, or
# text changed - the Title/Subtitle placeholders stay empty and
# ppt/slides/slide.xml is byte-for-byte identical before/after.
#
# The webextension part (ppt/slides/udata/data.xml) and the raw
# relationship-id values are internal to PowerPoint's task-pane-add-in
# plumbing and are not surfaced anywhere in the public Shape/Slide/
# Presentation COM object model (there is no Shapes entry, Tags, or
# CustomXMLPart for it - confirmed against the full real PowerPoint
# automation surface). So there is no COM call that flips that GUID.
# To stay faithful to the rest of the deck (whose OOXML is otherwise
# 100% unchanged) this script intentionally performs no shape/content
# mutations, touching only the already-open presentation object.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
